# Auto-update price data: insert a new row for the latest date at the top
# of the data table (row 2), pushing all existing rows down by one, and
# keeping a row for the oldest date at the bottom (duplicating the last
# row's values, matching the source feed's fixed-value behavior).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2 (shifts rows 2..70 down to 3..71).
$ws.Rows.Item(2).Insert()

# Seed the new row with the same cell formatting/values as the row right
# below it (the former top row), so styles/types stay identical to the
# rest of the table; we'll overwrite the date afterwards.
$ws.Range("A3:D3").Copy($ws.Range("A2:D2"))

# Write the new date as literal text (not an auto-converted date serial)
# by building it via a scratch-cell formula and pasting only the value.
# This preserves the same "plain text / no cell style" shape as every
# other date cell in the sheet.
$scratch = $ws.Range("Z1")
$scratch.Formula = '="2026-01-29"'
$scratch.Copy()
$ws.Range("A2").PasteSpecial(-4163)
$scratch.Clear()

# Numeric columns are unchanged constants for the new row.
$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
